# "Generate Report for Handback"
# Row 6 on both the zh-cn and de-de sheets corresponds to the
# db74cb5f-5c6c-4252-91e4-b7c9155840cf.md file. A handback report was
# generated for it: fill in the "Latest Target File", "Latest Handback
# File", "Latest Handback DateTime" and "Error Detail" columns (J, K, L, R)
# which were previously blank/placeholder, and widen the Error Detail
# column so the long message is readable.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/850d296c50eff54dca06124b04515681cbb40868/e2e/db74cb5f-5c6c-4252-91e4-b7c9155840cf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/966d08049f214dba4efdf4dd920412853e94e2b4/e2e/db74cb5f-5c6c-4252-91e4-b7c9155840cf.md."

$sheetInfo = @(
    @{ Name = "zh-cn"; Lang = "zh-cn"; HandbackDate = "2017-02-17 08:40:19"; TargetOrg = "ol-test4-zhcn"; TargetHash = "0e4c17f552ea9b350a70553561717b80954ee82c" },
    @{ Name = "de-de"; Lang = "de-de"; HandbackDate = "2017-02-17 08:40:43"; TargetOrg = "ol-test4-dede"; TargetHash = "097ed12bb79fdf063bc88fd809af1b34ebbfe6b9" }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    $targetFileName = "db74cb5f-5c6c-4252-91e4-b7c9155840cf.md"
    $handbackFileName = "db74cb5f-5c6c-4252-91e4-b7c9155840cf.d8c3270d535368f03b776ade556e6fff3fd5e980." + $info.Lang + ".xlf"

    # K6: Latest Handback File
    $ws.Range("K6").Value = $handbackFileName

    # L6: Latest Handback DateTime (plain text value, matches the other
    # date cells in this column which are stored as text, not dates)
    $ws.Range("L6").Value = $info.HandbackDate

    # R6: Error Detail
    $ws.Range("R6").Value = $errorMessage

    # J6: Latest Target File, as a hyperlink (same pattern as the other
    # rows in this column: display text is the bare filename, target is
    # a GitHub blob URL).
    $targetUrl = "https://github.com/OpenLocalizationTestOrg/" + $info.TargetOrg + "/blob/" + $info.TargetHash + "/e2e/" + $targetFileName
    $ws.Hyperlinks.Add($ws.Range("J6"), $targetUrl, "", "", $targetFileName)

    # Match the blue-underline look used by the other hyperlink cells
    # (column A and the other J cells) in this table.
    $ws.Range("J6").Font.Underline = 2
    $ws.Range("J6").Font.Color = 15570276

    # Widen the Error Detail column (R / column 18) so the long message
    # is readable, matching the width already used for columns A and G.
    $ws.Columns.Item(18).ColumnWidth = 40 - 5 / 6
}
